$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 31 ("e501 Game Won" row),
# shifting the existing rows 31-35 down to 33-37.
$ws.Rows("31:32").Insert()

# Column A labels first ("e030" / "e031"), then column B bodies - this
# matches the order the new shared-string entries were originally added.
$ws.Range("A31").Value = "e030"
$ws.Range("A32").Value = "e031"

# New row 31 body: e030 "Enemy Strength Roll Entering Battle Board"
$ws.Range("B31").Value = @'
<Bold>e030 Enemy Strength Roll Entering Battle Board</Bold> 
<InlineUIContainer><Button Content='r4.53' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table for enemy strength: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer> 
'@

# New row 32 body: e031 "Resistance Table"
$ws.Range("B32").Value = @'
<Bold>e031 Resistance Table</Bold> 
<InlineUIContainer><Button Content='r4.54.5' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    
<LineBreak/><LineBreak/>
Roll 1D and consult the <InlineUIContainer><Button Content='Resistance' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
Table to determine if combat occurs in this area: <LineBreak/><LineBreak/>
Die Roll =  <InlineUIContainer><Image Name='DieRoll' Height='21' Width='21' > </Image></InlineUIContainer> 
'@

# Both new rows hold multi-line wrapped text similar in length to the
# existing "e503/e504" rows, so set their heights to match the final file.
$ws.Range("A31").RowHeight = 90
$ws.Range("A32").RowHeight = 90

# Update view state: move the active selection to the newly inserted
# B31 cell (previously B27).
$ws.Range("B31").Select()
